# "Cadastro de cliente construído"
#
# Changes applied to the "Product Backlog" sheet (first sheet):
#   1. Filter the Product_Backlog table so only rows whose "Status"
#      column equals "Em Aberto" stay visible (AutoFilter on the
#      table's 2nd column) -- this hides rows 2-21 (all "Finalizado")
#      except row 13, which is already "Em Aberto".
#   2. Underline the "Em Aberto" text in cell B28.
#   3. Update the sheet's selection/view to B28 (and drop the old
#      frozen top-left cell scroll position).
#   4. Set the page setup (paper size / orientation) for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Filter the table on the Status column (index 2, 1-based) to only
#    show "Em Aberto" rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Range.AutoFilter(2, @("Em Aberto"), 7)

# 2) Underline cell B28 ("Em Aberto").
$ws.Range("B28").Font.Underline = $true

# 3) Move the selection/active cell to B28 and refresh the view.
$ws.Activate()
$ws.Range("B28").Select()

# 4) Page setup for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
